$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: SKU 1007e2 -> 24008-5, name -> "Pendant Green-Blues Lamb" ---
# Give the SKU cell the same "grey, larger" look already used by the old
# H4 supplier cell, then bump its size/typeface the rest of the way.
$ws.Range("H4").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Font.Size = 14
$ws.Range("A3").Font.Name = "Helvetica Neue"

$ws.Range("A3").Value = "24008-5"
$ws.Range("C3").Value = "Pendant Green-Blues Lamb"
# B3 (category "Pendants") and L3 (width=2) stay the same.

# --- Row 4: SKU 31903 -> 24008-6, name -> "Pendant Green-Blues 2" ---
# Re-use the formatting we just built for A3.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)

$ws.Range("A4").Value = "24008-6"
$ws.Range("C4").Value = "Pendant Green-Blues 2"
# B4 (category "Earrings") stays the same.

# Remove the old supplier (H4) and color (P4) values that no longer apply.
$ws.Range("H4").Clear()
$ws.Range("P4").Clear()

# Add the new size (D4) and long description (F4) values.
$ws.Range("D4").Value = 12
$ws.Range("F4").Value = "Some long description"

# --- Remove the old rows 5-8 (1076-5 / 1076-6 / 1076-7 / 1091e records) ---
$ws.Rows("5:8").Delete()

# --- Formatting tweaks ---
# Column C becomes a fixed (non bestFit) width, a bit wider than before.
$ws.Columns("C").ColumnWidth = 24.6

# Data rows 3 & 4 grow slightly taller to match the new bigger SKU font.
$ws.Rows(3).RowHeight = 17
$ws.Rows(4).RowHeight = 17

# Restore the cursor to where the next row of data would be entered.
$ws.Range("B9").Select()
